$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grades")

# Update the input scores that changed (E7 and E8: 1 -> 0.9)
$ws.Range("E7").Value = 0.9
$ws.Range("E8").Value = 0.9

# Update the active cell / selection to E8
$ws.Activate()
$ws.Range("E8").Select()
